# edit.ps1 - apply the "raven.docx" revision described by the commit diff.
# Uses raw WordprocessingML injection (Range.InsertXML) for exact,
# surgical control over run/paragraph structure, plus Find/Replace where
# a simple text change suffices.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) First paragraph: "This is a Microsoft word document." gets two
#    trailing spaces appended, followed by three new red (C00000) runs
#    spelling out "(This is a change – Version for branch alternate)".
# ---------------------------------------------------------------------
$p1 = $d.Paragraphs(1)
$r1 = $p1.Range
$xml1 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:r><w:t xml:space="preserve">This is a Microsoft word document.  </w:t></w:r>' +
        '<w:r><w:rPr><w:color w:val="C00000"/></w:rPr><w:t>(This is a change – Ve</w:t></w:r>' +
        '<w:r><w:rPr><w:color w:val="C00000"/></w:rPr><w:t>rsion for branch alternate</w:t></w:r>' +
        '<w:r><w:rPr><w:color w:val="C00000"/></w:rPr><w:t>)</w:t></w:r>' +
        '</w:p>'
$r1.InsertXML($xml1)

# ---------------------------------------------------------------------
# 2) Replace the empty paragraph right after "It will be treated as a
#    binary file by Git." with a new (still empty) shaded paragraph
#    (light-grey F9F9F9 background, bold Calibri 202122 run props).
#    This paragraph sits right before "The Raven" heading paragraph.
# ---------------------------------------------------------------------
$p3 = $d.Paragraphs(3)
$r3 = $p3.Range
$xml3 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:pPr>' +
          '<w:shd w:val="clear" w:color="auto" w:fill="F9F9F9"/>' +
          '<w:rPr>' +
            '<w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/>' +
            '<w:b/>' +
            '<w:bCs/>' +
            '<w:color w:val="202122"/>' +
          '</w:rPr>' +
        '</w:pPr>' +
        '</w:p>'
$r3.InsertXML($xml3)

# ---------------------------------------------------------------------
# 3) The trailing "ank God almighty, we are free at last." paragraph
#    (NormalWeb style, pasted-from-web remnants) is cleared out to a
#    bare empty paragraph, right before the sectPr.
# ---------------------------------------------------------------------
$lastIndex = $d.Paragraphs.Count
$pLast = $d.Paragraphs($lastIndex)
$rLast = $pLast.Range
$xmlLast = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>'
$rLast.InsertXML($xmlLast)
